# Apply the Jan 9 2024 cryptos-list refresh: updated prices/volume deltas
# for rows 2-46, and a 3-coin reshuffle (rows 47-51: ordi/Algorand swap,
# Aave/MultiversX swap) with their own refreshed price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    # Force the literal text even when it parses as a number (e.g. "305.91"),
    # matching the workbook convention of storing these as inline/shared strings
    # rather than numeric cells. Reset the style afterwards so no stray
    # number-format style is left attached to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "46.862.31"
$ws.Range("E2").Value = "  +6.40%  "
$ws.Range("D3").Value = "2.308.95"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextCell $ws.Range("D5") "305.91"
$ws.Range("E5").Value = "  +1.97%  "
Set-TextCell $ws.Range("D6") "102.57"
$ws.Range("E6").Value = "  +12.95%  "
Set-TextCell $ws.Range("D7") "0.573"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextCell $ws.Range("D9") "0.530"
$ws.Range("E9").Value = "  +7.36%  "
Set-TextCell $ws.Range("D10") "37.29"
$ws.Range("E10").Value = "  +12.50%  "
$ws.Range("E11").Value = "  +2.61%  "
Set-TextCell $ws.Range("D12") "7.50"
$ws.Range("E12").Value = "  +7.43%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "2.660.80"
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "2.299.96"
$ws.Range("E15").Value = "  +3.60%  "
Set-TextCell $ws.Range("D16") "14.12"
$ws.Range("E16").Value = "  +4.35%  "
Set-TextCell $ws.Range("D17") "0.823"
$ws.Range("E17").Value = "  +5.48%  "
$ws.Range("D18").Value = "46.840.25"
$ws.Range("E18").Value = "  +6.59%  "
Set-TextCell $ws.Range("D19") "13.55"
$ws.Range("E19").Value = "  +20.40%  "
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  +5.04%  "
$ws.Range("E21").Value = "  +2.56%  "
Set-TextCell $ws.Range("D22") "67.05"
$ws.Range("E22").Value = "  +3.37%  "
Set-TextCell $ws.Range("D23") "251.00"
$ws.Range("E23").Value = "  +5.73%  "
Set-TextCell $ws.Range("D24") "2.96"
$ws.Range("E24").Value = "  +4.61%  "
Set-TextCell $ws.Range("D25") "1.97"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("E26").Value = "  -0.19%  "
Set-TextCell $ws.Range("D27") "44.38"
$ws.Range("E27").Value = "  +14.90%  "
Set-TextCell $ws.Range("D28") "2.29"
$ws.Range("E28").Value = "  +5.28%  "
Set-TextCell $ws.Range("D29") "10.02"
$ws.Range("E29").Value = "  +6.95%  "
Set-TextCell $ws.Range("D30") "20.27"
$ws.Range("E30").Value = "  +4.89%  "
Set-TextCell $ws.Range("D31") "2.90"
$ws.Range("E31").Value = "  +16.24%  "
Set-TextCell $ws.Range("D32") "5.81"
$ws.Range("E32").Value = "  +6.97%  "
Set-TextCell $ws.Range("D33") "147.89"
$ws.Range("E33").Value = "  -1.98%  "
Set-TextCell $ws.Range("D34") "0.0807"
$ws.Range("E34").Value = "  +7.26%  "
Set-TextCell $ws.Range("D35") "3.22"
$ws.Range("E35").Value = "  +12.00%  "
$ws.Range("E36").Value = "  +11.58%  "
$ws.Range("E37").Value = "  +3.23%  "
Set-TextCell $ws.Range("D38") "1.82"
$ws.Range("E38").Value = "  +6.43%  "
Set-TextCell $ws.Range("D39") "16.35"
$ws.Range("E39").Value = "  +22.54%  "
Set-TextCell $ws.Range("D40") "4.17"
$ws.Range("E40").Value = "  +14.95%  "
Set-TextCell $ws.Range("D41") "3.48"
$ws.Range("E41").Value = "  +7.74%  "
Set-TextCell $ws.Range("D42") "0.0307"
$ws.Range("E42").Value = "  +0.75%  "
Set-TextCell $ws.Range("D43") "2.01"
$ws.Range("E43").Value = "  +10.72%  "
Set-TextCell $ws.Range("D44") "0.998"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "1.853.32"
$ws.Range("E45").Value = "  +1.41%  "
Set-TextCell $ws.Range("D46") "89.31"
$ws.Range("E46").Value = "  +21.10%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws.Range("D47") "0.200"
$ws.Range("E47").Value = "  +10.35%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextCell $ws.Range("D48") "75.14"
$ws.Range("E48").Value = "  +11.48%  "
Set-TextCell $ws.Range("D49") "4.96"
$ws.Range("E49").Value = "  +11.53%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell $ws.Range("D50") "55.20"
$ws.Range("E50").Value = "  +8.12%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Range("D51") "97.29"
$ws.Range("E51").Value = "  +3.04%  "
